$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells in column D are plain-text numeric-looking strings (e.g.
# "26.202.09", "19.20") in the source workbook. Assigning them straight to
# .Value would let Excel auto-coerce single-dot values into real numbers
# (dropping trailing zeros, e.g. "19.20" -> 19.2). Prefixing with a leading
# apostrophe forces text entry; ClearFormats() then strips the resulting
# quotePrefix cell style so the cell's styling matches the untouched original
# (no `s` attribute), leaving only the text value changed.
$ws.Range('D2').Value = "'26.202.09"
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -0.60%  '
$ws.Range('D3').Value = "'1.589.00"
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -0.07%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = "'211.81"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.75%  '
$ws.Range('D6').Value = "'0.503"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.30%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('E8').Value = '  -0.15%  '
$ws.Range('E9').Value = '  -0.94%  '
$ws.Range('D10').Value = "'19.20"
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -1.63%  '
$ws.Range('D11').Value = "'0.0847"
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.29%  '
$ws.Range('D12').Value = "'1.812.68"
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -0.08%  '
$ws.Range('D13').Value = "'1.602.05"
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.71%  '
$ws.Range('E14').Value = '  -1.53%  '
$ws.Range('E15').Value = '  -0.91%  '
$ws.Range('D16').Value = "'63.82"
$ws.Range('D16').ClearFormats()
$ws.Range('D17').Value = "'26.231.90"
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -0.45%  '
$ws.Range('D18').Value = "'0.0₃0725"
$ws.Range('D18').ClearFormats()
$ws.Range('E19').Value = '  -0.53%  '
$ws.Range('D20').Value = "'214.03"
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +1.56%  '
$ws.Range('E22').Value = '  -0.59%  '
$ws.Range('E23').Value = '  +0.59%  '
$ws.Range('D24').Value = "'2.11"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -1.65%  '
$ws.Range('D25').Value = "'144.30"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.57%  '
$ws.Range('E26').Value = '  -0.07%  '
$ws.Range('E27').Value = '  -0.98%  '
$ws.Range('E28').Value = '  -0.99%  '
$ws.Range('D29').Value = "'15.10"
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.93%  '
$ws.Range('E30').Value = '  -2.09%  '
$ws.Range('E31').Value = '  +0.30%  '
$ws.Range('E32').Value = '  -1.06%  '
$ws.Range('D33').Value = "'1.418.66"
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +8.57%  '
$ws.Range('D34').Value = "'2.96"
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -1.12%  '
$ws.Range('E35').Value = '  -0.60%  '
$ws.Range('D36').Value = "'0.588"
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -4.12%  '
$ws.Range('E38').Value = '  -1.38%  '
$ws.Range('D39').Value = "'5.92"
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +5.19%  '
$ws.Range('D40').Value = "'0.821"
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +1.36%  '
$ws.Range('E41').Value = '  -0.14%  '
$ws.Range('D42').Value = "'0.940"
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -14.03%  '
$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D43').Value = "'2.13"
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +0.10%  '
$ws.Range('B44').Value = 'TrustWalletToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D44').Value = "'0.764"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.28%  '
$ws.Range('D45').Value = "'1.724.03"
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.24%  '
$ws.Range('D46').Value = "'61.15"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -2.25%  '
$ws.Range('D47').Value = "'85.90"
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -2.32%  '
$ws.Range('E48').Value = '  -1.78%  '
$ws.Range('E49').Value = '  -0.50%  '
$ws.Range('E50').Value = '  -0.51%  '
$ws.Range('D51').Value = "'0.0968"
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -1.31%  '
